$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing column B validation values (rows 2-127) for the rerun on the new dataset
$ws.Cells.Item(2, 2).Value = 1.19685
$ws.Cells.Item(3, 2).Value = 1.16442
$ws.Cells.Item(4, 2).Value = 1.15956
$ws.Cells.Item(5, 2).Value = 1.17902
$ws.Cells.Item(6, 2).Value = 1.19685
$ws.Cells.Item(7, 2).Value = 1.16118
$ws.Cells.Item(8, 2).Value = 1.16767
$ws.Cells.Item(9, 2).Value = 1.15956
$ws.Cells.Item(10, 2).Value = 1.16605
$ws.Cells.Item(11, 2).Value = 1.17577
$ws.Cells.Item(12, 2).Value = 1.17902
$ws.Cells.Item(13, 2).Value = 1.18713
$ws.Cells.Item(14, 2).Value = 1.20658
$ws.Cells.Item(15, 2).Value = 1.22604
$ws.Cells.Item(16, 2).Value = 1.22928
$ws.Cells.Item(17, 2).Value = 1.21793
$ws.Cells.Item(18, 2).Value = 1.21631
$ws.Cells.Item(19, 2).Value = 1.22928
$ws.Cells.Item(20, 2).Value = 1.25847
$ws.Cells.Item(21, 2).Value = 1.23253
$ws.Cells.Item(22, 2).Value = 1.24874
$ws.Cells.Item(23, 2).Value = 1.24874
$ws.Cells.Item(24, 2).Value = 1.21307
$ws.Cells.Item(25, 2).Value = 1.2228
$ws.Cells.Item(26, 2).Value = 1.22442
$ws.Cells.Item(27, 2).Value = 1.23091
$ws.Cells.Item(28, 2).Value = 1.2228
$ws.Cells.Item(29, 2).Value = 1.2455
$ws.Cells.Item(30, 2).Value = 1.26496
$ws.Cells.Item(31, 2).Value = 1.26009
$ws.Cells.Item(32, 2).Value = 1.25523
$ws.Cells.Item(34, 2).Value = 1.24388
$ws.Cells.Item(35, 2).Value = 1.25036
$ws.Cells.Item(36, 2).Value = 1.26496
$ws.Cells.Item(37, 2).Value = 1.2682
$ws.Cells.Item(38, 2).Value = 1.26658
$ws.Cells.Item(39, 2).Value = 1.25523
$ws.Cells.Item(40, 2).Value = 1.22442
$ws.Cells.Item(41, 2).Value = 1.22928
$ws.Cells.Item(42, 2).Value = 1.22442
$ws.Cells.Item(43, 2).Value = 1.23577
$ws.Cells.Item(44, 2).Value = 1.24064
$ws.Cells.Item(45, 2).Value = 1.24064
$ws.Cells.Item(46, 2).Value = 1.23901
$ws.Cells.Item(47, 2).Value = 1.25036
$ws.Cells.Item(48, 2).Value = 1.25199
$ws.Cells.Item(49, 2).Value = 1.23901
$ws.Cells.Item(50, 2).Value = 1.22928
$ws.Cells.Item(51, 2).Value = 1.23091
$ws.Cells.Item(52, 2).Value = 1.24388
$ws.Cells.Item(53, 2).Value = 1.24388
$ws.Cells.Item(54, 2).Value = 1.23901
$ws.Cells.Item(55, 2).Value = 1.24064
$ws.Cells.Item(56, 2).Value = 1.30874
$ws.Cells.Item(57, 2).Value = 1.34279
$ws.Cells.Item(58, 2).Value = 1.36712
$ws.Cells.Item(59, 2).Value = 1.37684
$ws.Cells.Item(60, 2).Value = 1.38982
$ws.Cells.Item(61, 2).Value = 1.37522
$ws.Cells.Item(62, 2).Value = 1.39306
$ws.Cells.Item(63, 2).Value = 1.39144
$ws.Cells.Item(64, 2).Value = 1.38009
$ws.Cells.Item(65, 2).Value = 1.40603
$ws.Cells.Item(66, 2).Value = 1.38495
$ws.Cells.Item(67, 2).Value = 1.3736
$ws.Cells.Item(68, 2).Value = 1.44008
$ws.Cells.Item(69, 2).Value = 1.43846
$ws.Cells.Item(70, 2).Value = 1.44819
$ws.Cells.Item(71, 2).Value = 1.4563
$ws.Cells.Item(72, 2).Value = 1.48549
$ws.Cells.Item(73, 2).Value = 1.49522
$ws.Cells.Item(74, 2).Value = 1.48387
$ws.Cells.Item(75, 2).Value = 1.45792
$ws.Cells.Item(76, 2).Value = 1.49359
$ws.Cells.Item(77, 2).Value = 1.47251
$ws.Cells.Item(78, 2).Value = 1.47251
$ws.Cells.Item(79, 2).Value = 1.50332
$ws.Cells.Item(80, 2).Value = 1.58602
$ws.Cells.Item(81, 2).Value = 1.58602
$ws.Cells.Item(82, 2).Value = 1.56981
$ws.Cells.Item(83, 2).Value = 1.52116
$ws.Cells.Item(84, 2).Value = 1.50495
$ws.Cells.Item(85, 2).Value = 1.54548
$ws.Cells.Item(86, 2).Value = 1.57305
$ws.Cells.Item(87, 2).Value = 1.55197
$ws.Cells.Item(88, 2).Value = 1.53738
$ws.Cells.Item(89, 2).Value = 1.55359
$ws.Cells.Item(90, 2).Value = 1.55683
$ws.Cells.Item(91, 2).Value = 1.58764
$ws.Cells.Item(92, 2).Value = 1.57791
$ws.Cells.Item(93, 2).Value = 1.58278
$ws.Cells.Item(94, 2).Value = 1.58602
$ws.Cells.Item(95, 2).Value = 1.55521
$ws.Cells.Item(96, 2).Value = 1.53738
$ws.Cells.Item(97, 2).Value = 1.53251
$ws.Cells.Item(98, 2).Value = 1.51467
$ws.Cells.Item(99, 2).Value = 1.52116
$ws.Cells.Item(100, 2).Value = 1.53413
$ws.Cells.Item(101, 2).Value = 1.55197
$ws.Cells.Item(102, 2).Value = 1.5617
$ws.Cells.Item(103, 2).Value = 1.49846
$ws.Cells.Item(104, 2).Value = 1.52116
$ws.Cells.Item(105, 2).Value = 1.5017
$ws.Cells.Item(106, 2).Value = 1.4563
$ws.Cells.Item(107, 2).Value = 1.47738
$ws.Cells.Item(108, 2).Value = 1.50981
$ws.Cells.Item(109, 2).Value = 1.63629
$ws.Cells.Item(110, 2).Value = 1.68656
$ws.Cells.Item(111, 2).Value = 1.68656
$ws.Cells.Item(112, 2).Value = 1.75304
$ws.Cells.Item(113, 2).Value = 1.73034
$ws.Cells.Item(114, 2).Value = 1.7579
$ws.Cells.Item(115, 2).Value = 1.8179
$ws.Cells.Item(116, 2).Value = 1.81304
$ws.Cells.Item(117, 2).Value = 1.80169
$ws.Cells.Item(118, 2).Value = 1.81304
$ws.Cells.Item(119, 2).Value = 1.79194
$ws.Cells.Item(120, 2).Value = 1.80493
$ws.Cells.Item(121, 2).Value = 1.8633
$ws.Cells.Item(122, 2).Value = 1.88438
$ws.Cells.Item(123, 2).Value = 1.86979
$ws.Cells.Item(124, 2).Value = 1.84547
$ws.Cells.Item(125, 2).Value = 1.89411
$ws.Cells.Item(126, 2).Value = 1.92654
$ws.Cells.Item(127, 2).Value = 1.91844

# Extend the sheet with 6 new rows (128-133) for indices 126-131,
# copying column A formatting from the row above so the style index matches
$ws.Range("A122:A127").Copy()
$ws.Range("A128:A133").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(128, 1).Value = 126
$ws.Cells.Item(128, 2).Value = 1.82763
$ws.Cells.Item(129, 1).Value = 127
$ws.Cells.Item(129, 2).Value = 1.84222
$ws.Cells.Item(130, 1).Value = 128
$ws.Cells.Item(130, 2).Value = 1.91033
$ws.Cells.Item(131, 1).Value = 129
$ws.Cells.Item(131, 2).Value = 1.88925
$ws.Cells.Item(132, 1).Value = 130
$ws.Cells.Item(132, 2).Value = 1.8779
$ws.Cells.Item(133, 1).Value = 131
$ws.Cells.Item(133, 2).Value = 1.87141
